# Refresh crypto price/volume data (scheduled GitHub Actions update).
# Note: price cells that look like plain numbers are written with a leading
# apostrophe and then restyled to "Normal" so Excel keeps them as text
# (matching the source data, e.g. "87.20" instead of being coerced to 87.2),
# without leaving behind a lingering text-format style on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.912.06"
$ws.Range("E2").Value = "  +5.51%  "
$ws.Range("D3").Value = "3.650.56"
$ws.Range("E3").Value = "  +16.59%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "'594.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.54%  "
$ws.Range("D6").Value = "'181.88"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.55%  "
$ws.Range("D7").Value = "3.648.06"
$ws.Range("E7").Value = "  +16.60%  "
$ws.Range("E8").Value = "  +0.21%  "
$ws.Range("D9").Value = "'0.534"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.17%  "
$ws.Range("D10").Value = "'0.161"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.44%  "
$ws.Range("D11").Value = "'6.57"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.31%  "
$ws.Range("D12").Value = "'0.495"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.90%  "
$ws.Range("D13").Value = "'40.42"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +10.87%  "
$ws.Range("D14").Value = "'0.0000252"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.81%  "
$ws.Range("D15").Value = "4.260.55"
$ws.Range("E15").Value = "  +16.70%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "70.921.59"
$ws.Range("E16").Value = "  +5.67%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.662.76"
$ws.Range("E17").Value = "  +17.14%  "
$ws.Range("E18").Value = "  +0.83%  "
$ws.Range("D19").Value = "'7.44"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.80%  "
$ws.Range("D20").Value = "'16.88"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.46%  "
$ws.Range("D21").Value = "'512.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.04%  "
$ws.Range("D22").Value = "'9.11"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +15.65%  "
$ws.Range("D23").Value = "'0.737"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.89%  "
$ws.Range("D24").Value = "'87.20"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.66%  "
$ws.Range("E25").Value = "  +8.08%  "
$ws.Range("D26").Value = "'13.41"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.23%  "
$ws.Range("D27").Value = "'10.87"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.85%  "
$ws.Range("E28").Value = "  -0.13%  "
$ws.Range("D29").Value = "'2.51"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +8.68%  "
$ws.Range("D30").Value = "'8.11"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.58%  "
$ws.Range("B31").Value = "PEPE"
$ws.Range("C31").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D31").Value = "'0.0000111"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +16.99%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "'2.76"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.61%  "
$ws.Range("D33").Value = "'31.37"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +10.97%  "
$ws.Range("D34").Value = "'0.115"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.22%  "
$ws.Range("E35").Value = "  +0.23%  "
$ws.Range("D36").Value = "'6.07"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.66%  "
$ws.Range("E37").Value = "  +6.46%  "
$ws.Range("D38").Value = "'0.343"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +9.72%  "
$ws.Range("D39").Value = "'2.14"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.37%  "
$ws.Range("D40").Value = "'50.96"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.98%  "
$ws.Range("D41").Value = "'0.129"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.06%  "
$ws.Range("D42").Value = "'45.19"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.16%  "
$ws.Range("D43").Value = "3.126.81"
$ws.Range("E43").Value = "  +11.16%  "
$ws.Range("D44").Value = "'8.78"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.98%  "
$ws.Range("D45").Value = "'411.45"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +9.46%  "
$ws.Range("E46").Value = "  +2.58%  "
$ws.Range("D47").Value = "'0.0367"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.78%  "
$ws.Range("D48").Value = "'28.12"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +13.18%  "
$ws.Range("D49").Value = "'137.39"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.87%  "
$ws.Range("E51").Value = "  +10.59%  "
